# Update weights and eggs - append new tracking rows to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily log entries (dates stored as Excel serial numbers, formatted as dates
# via the same style already used in column A for the existing rows)
$ws.Range("A19").Value = 44233
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0

$ws.Range("A20").Value = 44234
$ws.Range("B20").Value = 24
$ws.Range("C20").Value = 15

$ws.Range("A21").Value = 44235
$ws.Range("B21").Value = 8
$ws.Range("C21").Value = 0

# Match the date number format already applied to the rest of column A
$ws.Range("A18").Copy()
$ws.Range("A19:A21").PasteSpecial(-4122)

# Reflect the final selection state from the edit
$ws.Range("H22").Select()
